$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (not auto-converted to a number),
# matching the inlineStr cell type used throughout this sheet,
# then restore the default "Normal" style so no stray number format sticks.
function Set-TextValue {
    param($Cell, $Text)
    $r = $ws.Range($Cell)
    $r.NumberFormat = "@"
    $r.Value = $Text
    $r.Style = "Normal"
}

# Row 2
Set-TextValue "D2" '44.554.31'
$ws.Range("E2").Value = '  +1.30%  '

# Row 3
Set-TextValue "D3" '2.242.83'
$ws.Range("E3").Value = '  +0.22%  '

# Row 4
$ws.Range("E4").Value = '  +0.72%  '

# Row 5
Set-TextValue "D5" '306.66'
$ws.Range("E5").Value = '  +0.02%  '

# Row 6
Set-TextValue "D6" '95.19'
$ws.Range("E6").Value = '  +0.18%  '

# Row 7
Set-TextValue "D7" '0.571'
$ws.Range("E7").Value = '  +0.34%  '

# Row 8
$ws.Range("E8").Value = '  +0.15%  '

# Row 9
$ws.Range("E9").Value = '  +0.29%  '

# Row 10
$ws.Range("E10").Value = '  +0.40%  '

# Row 11
$ws.Range("E11").Value = '  -0.43%  '

# Row 12
Set-TextValue "D12" '7.23'
$ws.Range("E12").Value = '  +0.15%  '

# Row 13
$ws.Range("E13").Value = '  +0.17%  '

# Row 14
Set-TextValue "D14" '2.281.24'
$ws.Range("E14").Value = '  +2.09%  '

# Row 15
Set-TextValue "D15" '0.835'
$ws.Range("E15").Value = '  +1.16%  '

# Row 16
Set-TextValue "D16" '13.60'
$ws.Range("E16").Value = '  +0.03%  '

# Row 17
Set-TextValue "D17" '44.284.77'
$ws.Range("E17").Value = '  +0.93%  '

# Row 18
$ws.Range("E18").Value = '  -0.82%  '

# Row 19
Set-TextValue "D19" '6.33'
$ws.Range("E19").Value = '  +1.14%  '

# Row 20
Set-TextValue "D20" '11.98'
$ws.Range("E20").Value = '  -1.20%  '

# Row 21
Set-TextValue "D21" '65.60'
$ws.Range("E21").Value = '  +0.80%  '

# Row 22
Set-TextValue "D22" '237.73'
$ws.Range("E22").Value = '  +0.50%  '

# Row 23
$ws.Range("E23").Value = '  +0.84%  '

# Row 24
$ws.Range("E24").Value = '  +1.07%  '

# Row 25
$ws.Range("E25").Value = '  -0.03%  '

# Row 26
Set-TextValue "D26" '2.23'
$ws.Range("E26").Value = '  +2.98%  '

# Row 27
Set-TextValue "D27" '37.89'
$ws.Range("E27").Value = '  +1.53%  '

# Row 28
Set-TextValue "D28" '9.80'
$ws.Range("E28").Value = '  -1.54%  '

# Row 29
Set-TextValue "D29" '5.99'
$ws.Range("E29").Value = '  +0.04%  '

# Row 30
Set-TextValue "D30" '19.95'
$ws.Range("E30").Value = '  +0.21%  '

# Row 31
Set-TextValue "D31" '153.01'
$ws.Range("E31").Value = '  +0.09%  '

# Row 32
$ws.Range("E32").Value = '  -0.53%  '

# Row 33
$ws.Range("E33").Value = '  +3.08%  '

# Row 34
Set-TextValue "D34" '3.04'
$ws.Range("E34").Value = '  -6.81%  '

# Row 35
Set-TextValue "D35" '0.111'
$ws.Range("E35").Value = '  +2.11%  '

# Row 36
$ws.Range("E36").Value = '  +0.39%  '

# Row 37
$ws.Range("E37").Value = '  +2.49%  '

# Row 38
Set-TextValue "D38" '14.99'
$ws.Range("E38").Value = '  -1.15%  '

# Row 39
Set-TextValue "D39" '3.39'
$ws.Range("E39").Value = '  +1.24%  '

# Row 40
$ws.Range("E40").Value = '  -1.88%  '

# Row 41
$ws.Range("E41").Value = '  +0.08%  '

# Row 43
Set-TextValue "D43" '1.793.16'
$ws.Range("E43").Value = '  +3.89%  '

# Row 44
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue "D44" '0.192'
$ws.Range("E44").Value = '  +2.28%  '

# Row 45
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue "D45" '1.68'
$ws.Range("E45").Value = '  +11.82%  '

# Row 46
Set-TextValue "D46" '79.07'
$ws.Range("E46").Value = '  -7.34%  '

# Row 47
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue "D47" '99.01'
$ws.Range("E47").Value = '  -1.04%  '

# Row 48
$ws.Range("B48").Value = 'THORChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue "D48" '4.92'
$ws.Range("E48").Value = '  -0.24%  '

# Row 49
$ws.Range("B49").Value = 'ordi'
$ws.Range("C49").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
Set-TextValue "D49" '70.20'
$ws.Range("E49").Value = '  +1.30%  '

# Row 50
Set-TextValue "D50" '8.13'
$ws.Range("E50").Value = '  +0.73%  '

# Row 51
Set-TextValue "D51" '54.61'
$ws.Range("E51").Value = '  +0.77%  '
